$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9358481262327414
$ws.Range("C2").Value = 0.05627632771751921

$ws.Range("B3").Value = 0.8342702169625247
$ws.Range("C3").Value = 0.07848526561490361

$ws.Range("B4").Value = 0.9253944773175542
$ws.Range("C4").Value = 0.1004054357682348

$ws.Range("B5").Value = 0.7519723865877712
$ws.Range("C5").Value = 0.07430820862156795

$ws.Range("B6").Value = 0.8819156804733728
$ws.Range("C6").Value = 0.0792477122076556
